$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '44.000.49'
$ws.Range('E2').Value = '  -1.06%  '

# Row 3
$ws.Range('D3').Value = '2.241.97'
$ws.Range('E3').Value = '  -1.66%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.55'
$ws.Range('E5').Value = '  -1.71%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.10'
$ws.Range('E6').Value = '  -6.61%  '

# Row 7
$ws.Range('E7').Value = '  -3.10%  '

# Row 8
$ws.Range('E8').Value = '  +0.11%  '

# Row 9
$ws.Range('E9').Value = '  -6.93%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.21'
$ws.Range('E10').Value = '  -6.67%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0821'
$ws.Range('E11').Value = '  -2.75%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.35'
$ws.Range('E12').Value = '  -6.98%  '

# Row 13
$ws.Range('E13').Value = '  -2.86%  '

# Row 14
$ws.Range('D14').Value = '2.582.71'
$ws.Range('E14').Value = '  -1.76%  '

# Row 15
$ws.Range('E15').Value = '  -4.90%  '

# Row 16
$ws.Range('D16').Value = '2.245.56'
$ws.Range('E16').Value = '  -1.51%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.96'
$ws.Range('E17').Value = '  -4.57%  '

# Row 18
$ws.Range('D18').Value = '43.881.41'
$ws.Range('E18').Value = '  -1.05%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.14'
$ws.Range('E19').Value = '  -6.72%  '

# Row 20
$ws.Range('E20').Value = '  -2.88%  '

# Row 21
$ws.Range('E21').Value = '  -3.15%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '65.54'
$ws.Range('E22').Value = '  -1.36%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.84'
$ws.Range('E23').Value = '  -1.06%  '

# Row 24
$ws.Range('E24').Value = '  -7.29%  '

# Row 25
$ws.Range('E25').Value = '  -8.55%  '

# Row 26
$ws.Range('E26').Value = '  +0.19%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.13'
$ws.Range('E27').Value = '  -0.75%  '

# Row 28
$ws.Range('E28').Value = '  -4.54%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.40'
$ws.Range('E29').Value = '  -5.15%  '

# Row 30
$ws.Range('E30').Value = '  -8.42%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.07'
$ws.Range('E31').Value = '  -2.82%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '155.77'
$ws.Range('E32').Value = '  -4.87%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0839'
$ws.Range('E33').Value = '  -5.31%  '

# Row 34
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.66'
$ws.Range('E34').Value = '  -3.22%  '

# Row 35
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.28'
$ws.Range('E35').Value = '  +2.65%  '

# Row 36
$ws.Range('E36').Value = '  -7.55%  '

# Row 37
$ws.Range('E37').Value = '  -7.89%  '

# Row 38
$ws.Range('E38').Value = '  -3.05%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '15.42'
$ws.Range('E39').Value = '  -1.06%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.52'
$ws.Range('E40').Value = '  -11.55%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.99'
$ws.Range('E41').Value = '  -10.69%  '

# Row 42
$ws.Range('E42').Value = '  -6.54%  '

# Row 43
$ws.Range('E43').Value = '  +0.12%  '

# Row 44
$ws.Range('D44').Value = '1.698.39'
$ws.Range('E44').Value = '  -4.27%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '82.88'
$ws.Range('E45').Value = '  -4.56%  '

# Row 46
$ws.Range('E46').Value = '  -6.41%  '

# Row 47
$ws.Range('E47').Value = '  -6.16%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '101.71'
$ws.Range('E48').Value = '  -2.62%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '71.04'
$ws.Range('E49').Value = '  -5.03%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '56.31'
$ws.Range('E50').Value = '  -6.71%  '

# Row 51
$ws.Range('E51').Value = '  -5.87%  '
